$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($ws, $rowNum, $values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($rowNum, $i + 1).Value = $values[$i]
    }
}

# Insert a new row at position 4 (pushes old rows 4-7 down to 5-8),
# so the new "E3" category data lands on row 4. The previously-existing
# rows (old E4/E5/E6/E7 @ rows 4-7) move intact to rows 5-8, but their
# category index in column A needs to be bumped by one since a new
# category was inserted ahead of them.
$ws.Rows.Item(4).Insert()

# Row 2 - E1 (values recomputed after the new experiment run)
Set-RowValues $ws 2 @(0, "E1", 11, 480, 11, 0, 0, 0, 0.02291666666666667, 0.7857142857142857, 0.04453441295546558, 0.02291666666666667, 0.7857142857142857, 0.04453441295546558)

# Row 3 - E2 (values recomputed after the new experiment run)
Set-RowValues $ws 3 @(1, "E2", 14, 932, 14, 0, 0, 0, 0.01502145922746781, 1, 0.02959830866807611, 0.01502145922746781, 1, 0.02959830866807611)

# Row 4 - E3 (brand new category)
Set-RowValues $ws 4 @(2, "E3", 11, 442, 11, 0, 0, 0, 0.0248868778280543, 0.7857142857142857, 0.04824561403508772, 0.0248868778280543, 0.7857142857142857, 0.04824561403508772)

# Rows 5-8 kept their old data intact via the row insert/shift, but column A
# (the category index) needs to be bumped by one to account for the new row.
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(8, 1).Value = 6

# Give the new row's category-index cell (A4) the same bold/bordered/centered
# formatting used by the rest of column A.
$src = $ws.Cells.Item(5, 1)
$dst = $ws.Cells.Item(4, 1)
$dst.Font.Bold = $src.Font.Bold
$dst.HorizontalAlignment = $src.HorizontalAlignment
$dst.VerticalAlignment = $src.VerticalAlignment
$dst.Borders.LineStyle = $src.Borders.LineStyle
